$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Best effort: make any new comment's author read "Author".
$excel.UserName = "Author"

# --- "Reporting Pending" row (row 5) ------------------------------------
# D5: "4 Dayes" -> "6 Dayes"
$ws.Range("D5").Value = "6 Dayes"

# New cells E5:G5 noting the still-open issue and its date window
$ws.Range("E5").Value = "Issue Pending"
$ws.Range("F5").Value = "20-09-21"
$ws.Range("G5").Value = "21-09-21"

# --- "Program Wise Delivery Challan Entry" row (row 7) -------------------
# C7: "Pending" -> "ok"
$ws.Range("C7").Value = "ok"

# --- Note on D5 explaining the date range ---------------------------------
$comment = $ws.Range("D5").AddComment("Author:" + [char]10 + "16-09-21 to 21-09-21")
$comment.Author = "Author"

# --- Match the author's final selection -----------------------------------
$ws.Range("D16").Select()
